$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Status column (C) for all data rows from "AVAILABLE" to "Available"
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 3).Value = "Available"
}

# Move the active cell selection to G9 to match the saved view state
$ws.Range("G9").Select()
